# SphereTests now running in terminal
# Applies the cell-value and selection changes described by the target diff.

$wb = $excel.ActiveWorkbook

$wsMain  = $wb.Worksheets.Item("MAIN Config.")
$wsComp  = $wb.Worksheets.Item("Computational benchmarks")
$wsExp   = $wb.Worksheets.Item("Experimental benchmarks")
$wsLib   = $wb.Worksheets.Item("Libraries")

# --- MAIN Config. sheet -----------------------------------------------
# MPI tasks: 32 -> 8
$wsMain.Range("B11").Value = 8

# --- Computational benchmarks sheet ------------------------------------
# Sphere Leakage Test row: NPS 10,000,000 -> 1,000,000 ; add NPS cut-off value 10
$wsComp.Range("I4").Value = 1000000
$wsComp.Range("J4").Value = 10

# Rows 5-7, column D ("OnlyInput"): true -> false.
# Use copy / paste-values from a neighbouring cell that already holds the
# shared string "false" with the identical style, so the written cell keeps
# its original text type + formatting instead of being auto-typed to a
# native boolean by plain value assignment.
$wsComp.Range("C5").Copy()
$wsComp.Range("D5").PasteSpecial(-4163)
$wsComp.Range("C6").Copy()
$wsComp.Range("D6").PasteSpecial(-4163)
$wsComp.Range("C7").Copy()
$wsComp.Range("D7").PasteSpecial(-4163)

# --- Experimental benchmarks sheet -------------------------------------
$wsExp.Range("C4").Copy()
$wsExp.Range("D4").PasteSpecial(-4163)

# --- Update the saved cursor/selection per sheet -----------------------
# Selecting a range activates its sheet, which would disturb which tab is
# marked active/selected, so restore the original active sheet
# ("Computational benchmarks") once all selections have been applied.
$wsMain.Range("B12").Select()
$wsComp.Range("C4").Select()
$wsExp.Range("E18").Select()
$wsLib.Range("D20").Select()

$wsComp.Activate()
